$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 2).Value = 23.97304912817162
$ws.Cells.Item(2, 3).Value = 9.467021466637208
$ws.Cells.Item(2, 4).Value = 7.496262721895841
$ws.Cells.Item(2, 6).Value = 40.33229155271498
$ws.Cells.Item(2, 7).Value = 3.693699081453545
$ws.Cells.Item(2, 12).Value = 10.8179865546421
$ws.Cells.Item(2, 14).Value = 19.86063773198385

# Row 3
$ws.Cells.Item(3, 2).Value = 23.44570026549173
$ws.Cells.Item(3, 3).Value = 8.864468073970492
$ws.Cells.Item(3, 4).Value = 7.518017032785777
$ws.Cells.Item(3, 6).Value = 39.93198758602487
$ws.Cells.Item(3, 7).Value = 3.698413125267562
$ws.Cells.Item(3, 12).Value = 10.79925974713658
$ws.Cells.Item(3, 14).Value = 19.91694477518601

# Row 4
$ws.Cells.Item(4, 2).Value = 23.12399964270126
$ws.Cells.Item(4, 3).Value = 8.47222758776892
$ws.Cells.Item(4, 4).Value = 7.53260962496968
$ws.Cells.Item(4, 6).Value = 39.69704444233493
$ws.Cells.Item(4, 7).Value = 3.701454010776515
$ws.Cells.Item(4, 12).Value = 10.79019025838891
$ws.Cells.Item(4, 14).Value = 19.95349934376373

# Row 5
$ws.Cells.Item(5, 2).Value = 22.99363543139072
$ws.Cells.Item(5, 3).Value = 8.3067846572525
$ws.Cells.Item(5, 4).Value = 7.538865535986485
$ws.Cells.Item(5, 6).Value = 39.60411711564352
$ws.Cells.Item(5, 7).Value = 3.702730179429174
$ws.Cells.Item(5, 12).Value = 10.78710641214416
$ws.Cells.Item(5, 14).Value = 19.96889349210754

# Row 6
$ws.Cells.Item(6, 2).Value = 22.97203877191087
$ws.Cells.Item(6, 3).Value = 8.278974681018667
$ws.Cells.Item(6, 4).Value = 7.539922968279894
$ws.Cells.Item(6, 6).Value = 39.58885896247192
$ws.Cells.Item(6, 7).Value = 3.702944324602085
$ws.Cells.Item(6, 12).Value = 10.78663133567222
$ws.Cells.Item(6, 14).Value = 19.9714797351258

# Row 7
$ws.Cells.Item(7, 2).Value = 23.1222382709865
$ws.Cells.Item(7, 3).Value = 8.4700190473971
$ws.Cells.Item(7, 4).Value = 7.532692743718545
$ws.Cells.Item(7, 6).Value = 39.69577968865197
$ws.Cells.Item(7, 7).Value = 3.701471071678459
$ws.Cells.Item(7, 12).Value = 10.79014618906556
$ws.Cells.Item(7, 14).Value = 19.95370493938239

# Row 8
$ws.Cells.Item(8, 2).Value = 23.79090658339463
$ws.Cells.Item(8, 3).Value = 9.263881018271753
$ws.Cells.Item(8, 4).Value = 7.503506406108061
$ws.Cells.Item(8, 6).Value = 40.19206992041448
$ws.Cells.Item(8, 7).Value = 3.695294185200878
$ws.Cells.Item(8, 12).Value = 10.81102636384295
$ws.Cells.Item(8, 14).Value = 19.87964078120646

# Row 9
$ws.Cells.Item(9, 2).Value = 25.1099835082657
$ws.Cells.Item(9, 3).Value = 10.64420799602097
$ws.Cells.Item(9, 4).Value = 7.456134421448116
$ws.Cells.Item(9, 6).Value = 41.24722762333881
$ws.Cells.Item(9, 7).Value = 3.684336046973797
$ws.Cells.Item(9, 12).Value = 10.87116705772443
$ws.Cells.Item(9, 14).Value = 19.75014733185512

# Row 10
$ws.Cells.Item(10, 2).Value = 26.07254759404719
$ws.Cells.Item(10, 3).Value = 11.55135948009817
$ws.Cells.Item(10, 4).Value = 7.427426566682736
$ws.Cells.Item(10, 6).Value = 42.0665223343106
$ws.Cells.Item(10, 7).Value = 3.676978894844527
$ws.Cells.Item(10, 12).Value = 10.92692142192756
$ws.Cells.Item(10, 14).Value = 19.66463702423567

# Row 11
$ws.Cells.Item(11, 2).Value = 26.50678551183399
$ws.Cells.Item(11, 3).Value = 11.94098347360784
$ws.Cells.Item(11, 4).Value = 7.415708843662698
$ws.Cells.Item(11, 6).Value = 42.44747584984778
$ws.Cells.Item(11, 7).Value = 3.673780406908087
$ws.Cells.Item(11, 12).Value = 10.95475908317635
$ws.Cells.Item(11, 14).Value = 19.62783407345927

# Row 12
$ws.Cells.Item(12, 2).Value = 26.67052464133744
$ws.Cells.Item(12, 3).Value = 12.0852246965661
$ws.Cells.Item(12, 4).Value = 7.411466131908632
$ws.Cells.Item(12, 6).Value = 42.59280384132298
$ws.Cells.Item(12, 7).Value = 3.672590381499908
$ws.Cells.Item(12, 12).Value = 10.96565223861232
$ws.Cells.Item(12, 14).Value = 19.61419994206807

# Row 13
$ws.Cells.Item(13, 2).Value = 26.6352941619037
$ws.Cells.Item(13, 3).Value = 12.05430636428162
$ws.Cells.Item(13, 4).Value = 7.412371197877174
$ws.Cells.Item(13, 6).Value = 42.56145919130672
$ws.Cells.Item(13, 7).Value = 3.672845735579861
$ws.Cells.Item(13, 12).Value = 10.96329063585514
$ws.Cells.Item(13, 14).Value = 19.61712283384076

# Row 14
$ws.Cells.Item(14, 2).Value = 26.52027136335199
$ws.Cells.Item(14, 3).Value = 11.95291636467891
$ws.Cells.Item(14, 4).Value = 7.415355885126286
$ws.Cells.Item(14, 6).Value = 42.45941132128518
$ws.Cells.Item(14, 7).Value = 3.673682079361614
$ws.Cells.Item(14, 12).Value = 10.9556482481856
$ws.Cells.Item(14, 14).Value = 19.62670632027454

# Row 15
$ws.Cells.Item(15, 2).Value = 26.4497207790109
$ws.Cells.Item(15, 3).Value = 11.89038258006707
$ws.Cells.Item(15, 4).Value = 7.417209474710332
$ws.Cells.Item(15, 6).Value = 42.397039780408
$ws.Cells.Item(15, 7).Value = 3.674197116758911
$ws.Cells.Item(15, 12).Value = 10.95101272509898
$ws.Cells.Item(15, 14).Value = 19.63261588541996

# Row 16
$ws.Cells.Item(16, 2).Value = 26.04408078149434
$ws.Cells.Item(16, 3).Value = 11.52543400068017
$ws.Cells.Item(16, 4).Value = 7.428219481649198
$ws.Cells.Item(16, 6).Value = 42.04178217221967
$ws.Cells.Item(16, 7).Value = 3.677190894801831
$ws.Cells.Item(16, 12).Value = 10.92515161284138
$ws.Cells.Item(16, 14).Value = 19.66708444098382

# Row 17
$ws.Cells.Item(17, 2).Value = 25.79417468709196
$ws.Cells.Item(17, 3).Value = 11.29565652099801
$ws.Cells.Item(17, 4).Value = 7.435318526593529
$ws.Cells.Item(17, 6).Value = 41.82587364060009
$ws.Cells.Item(17, 7).Value = 3.679065357787964
$ws.Cells.Item(17, 12).Value = 10.90991764077001
$ws.Cells.Item(17, 14).Value = 19.6887672928054

# Row 18
$ws.Cells.Item(18, 2).Value = 25.65010269704149
$ws.Cells.Item(18, 3).Value = 11.1613247638255
$ws.Cells.Item(18, 4).Value = 7.439527891870481
$ws.Cells.Item(18, 6).Value = 41.70247266362723
$ws.Cells.Item(18, 7).Value = 3.680157468946268
$ws.Cells.Item(18, 12).Value = 10.90138867938163
$ws.Cells.Item(18, 14).Value = 19.70143593344727

# Row 19
$ws.Cells.Item(19, 2).Value = 25.60127076072299
$ws.Cells.Item(19, 3).Value = 11.11546929879845
$ws.Cells.Item(19, 4).Value = 7.440974733773997
$ws.Cells.Item(19, 6).Value = 41.66082940330243
$ws.Cells.Item(19, 7).Value = 3.680529643046678
$ws.Cells.Item(19, 12).Value = 10.89854109557588
$ws.Cells.Item(19, 14).Value = 19.70575916578008

# Row 20
$ws.Cells.Item(20, 2).Value = 25.82081333948166
$ws.Cells.Item(20, 3).Value = 11.3203412059814
$ws.Cells.Item(20, 4).Value = 7.434549750138082
$ws.Cells.Item(20, 6).Value = 41.8487771882834
$ws.Cells.Item(20, 7).Value = 3.678864373325421
$ws.Cells.Item(20, 12).Value = 10.91151521517435
$ws.Cells.Item(20, 14).Value = 19.68643869631001

# Row 21
$ws.Cells.Item(21, 2).Value = 26.5540766073261
$ws.Cells.Item(21, 3).Value = 11.982786544864
$ws.Cells.Item(21, 4).Value = 7.414473917540853
$ws.Cells.Item(21, 6).Value = 42.48935718055437
$ws.Cells.Item(21, 7).Value = 3.67343585153992
$ws.Cells.Item(21, 12).Value = 10.95788349539401
$ws.Cells.Item(21, 14).Value = 19.62388320471085

# Row 22
$ws.Cells.Item(22, 2).Value = 27.02916894698076
$ws.Cells.Item(22, 3).Value = 12.39649913478821
$ws.Cells.Item(22, 4).Value = 7.402487865393086
$ws.Cells.Item(22, 6).Value = 42.91418989592478
$ws.Cells.Item(22, 7).Value = 3.670011336555359
$ws.Cells.Item(22, 12).Value = 10.99023516626941
$ws.Cells.Item(22, 14).Value = 19.58476235529511

# Row 23
$ws.Cells.Item(23, 2).Value = 26.77603651136388
$ws.Cells.Item(23, 3).Value = 12.17744836089539
$ws.Cells.Item(23, 4).Value = 7.408780687206548
$ws.Cells.Item(23, 6).Value = 42.68692263113109
$ws.Cells.Item(23, 7).Value = 3.67182783053507
$ws.Cells.Item(23, 12).Value = 10.97278265862698
$ws.Cells.Item(23, 14).Value = 19.60548027362438

# Row 24
$ws.Cells.Item(24, 2).Value = 25.80877122997256
$ws.Cells.Item(24, 3).Value = 11.30918820563304
$ws.Cells.Item(24, 4).Value = 7.434896915311056
$ws.Cells.Item(24, 6).Value = 41.83842021881426
$ws.Cells.Item(24, 7).Value = 3.678955193357635
$ws.Cells.Item(24, 12).Value = 10.91079223740236
$ws.Cells.Item(24, 14).Value = 19.68749082292852

# Row 25
$ws.Cells.Item(25, 2).Value = 24.75350087524398
$ws.Cells.Item(25, 3).Value = 10.2896294118907
$ws.Cells.Item(25, 4).Value = 7.467885367379553
$ws.Cells.Item(25, 6).Value = 40.95359719294225
$ws.Cells.Item(25, 7).Value = 3.687177943534206
$ws.Cells.Item(25, 12).Value = 10.85285415853767
$ws.Cells.Item(25, 14).Value = 19.78349016409963
